$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-08-02 Friday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-08-03 Saturday", 2) | Out-Null
$d.Content.Find.Execute("511×7=3577", $false, $false, $false, $false, $false, $true, 1, $false, "105×4=420", 2) | Out-Null
$d.Content.Find.Execute("135×8=1080", $false, $false, $false, $false, $false, $true, 1, $false, "519×4=2076", 2) | Out-Null
$d.Content.Find.Execute("931×5=4655", $false, $false, $false, $false, $false, $true, 1, $false, "352×3=1056", 2) | Out-Null
$d.Content.Find.Execute("206×3=618", $false, $false, $false, $false, $false, $true, 1, $false, "152×6=912", 2) | Out-Null
$d.Content.Find.Execute("110×4=440", $false, $false, $false, $false, $false, $true, 1, $false, "844×2=1688", 2) | Out-Null
$d.Content.Find.Execute("754×6=4524", $false, $false, $false, $false, $false, $true, 1, $false, "615×3=1845", 2) | Out-Null
$d.Content.Find.Execute("849×3=2547", $false, $false, $false, $false, $false, $true, 1, $false, "352×3=1056", 2) | Out-Null
$d.Content.Find.Execute("341×7=2387", $false, $false, $false, $false, $false, $true, 1, $false, "720×8=5760", 2) | Out-Null
$d.Content.Find.Execute("370×6=2220", $false, $false, $false, $false, $false, $true, 1, $false, "747×4=2988", 2) | Out-Null
$d.Content.Find.Execute("951×8=7608", $false, $false, $false, $false, $false, $true, 1, $false, "611×4=2444", 2) | Out-Null
$d.Content.Find.Execute("752×5=3760", $false, $false, $false, $false, $false, $true, 1, $false, "567×3=1701", 2) | Out-Null
$d.Content.Find.Execute("246×6=1476", $false, $false, $false, $false, $false, $true, 1, $false, "818×8=6544", 2) | Out-Null
$d.Content.Find.Execute("404×8=3232", $false, $false, $false, $false, $false, $true, 1, $false, "232×2=464", 2) | Out-Null
$d.Content.Find.Execute("571×9=5139", $false, $false, $false, $false, $false, $true, 1, $false, "877×8=7016", 2) | Out-Null
$d.Content.Find.Execute("417×2=834", $false, $false, $false, $false, $false, $true, 1, $false, "198×7=1386", 2) | Out-Null
$d.Content.Find.Execute("920×4=3680", $false, $false, $false, $false, $false, $true, 1, $false, "677×4=2708", 2) | Out-Null
$d.Content.Find.Execute("885×2=1770", $false, $false, $false, $false, $false, $true, 1, $false, "370×3=1110", 2) | Out-Null
$d.Content.Find.Execute("166×4=664", $false, $false, $false, $false, $false, $true, 1, $false, "220×2=440", 2) | Out-Null
$d.Content.Find.Execute("453×5=2265", $false, $false, $false, $false, $false, $true, 1, $false, "493×2=986", 2) | Out-Null
$d.Content.Find.Execute("304×2=608", $false, $false, $false, $false, $false, $true, 1, $false, "162×3=486", 2) | Out-Null
$d.Content.Find.Execute("173×9=1557", $false, $false, $false, $false, $false, $true, 1, $false, "402×5=2010", 2) | Out-Null
$d.Content.Find.Execute("221×7=1547", $false, $false, $false, $false, $false, $true, 1, $false, "610×8=4880", 2) | Out-Null
$d.Content.Find.Execute("900×6=5400", $false, $false, $false, $false, $false, $true, 1, $false, "755×8=6040", 2) | Out-Null
$d.Content.Find.Execute("765×3=2295", $false, $false, $false, $false, $false, $true, 1, $false, "434×8=3472", 2) | Out-Null
$d.Content.Find.Execute("194×3=582", $false, $false, $false, $false, $false, $true, 1, $false, "374×5=1870", 2) | Out-Null
